$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.99216806615691489
$ws.Range("A3").Value = 0.84828481282431789
$ws.Range("B3").Value = 0.93754551010813525
$ws.Range("AB3").Value = 0.87525963095701487
$ws.Range("C4").Value = 0.9178394100096765
$ws.Range("E4").Value = 0.9522930783198239
$ws.Range("C5").Value = 0.70287659127588897
$ws.Range("D6").Value = 0.97818388928284428
$ws.Range("E6").Value = 0.64170888569924101
$ws.Range("H6").Value = 0.6766203607773601
$ws.Range("L6").Value = 0.78302347585209287
$ws.Range("V6").Value = 0.87037935590009974
$ws.Range("F7").Value = 0.71038930686333635
$ws.Range("H7").Value = 0.81946466074977531
$ws.Range("I7").Value = 0.86443892114094822
$ws.Range("S7").Value = 0.93225108114264299
$ws.Range("H9").Value = 0.71901235109475148
$ws.Range("J9").Value = 0.96557468939225699
$ws.Range("AU9").Value = 0.98248705987410245
$ws.Range("BJ9").Value = 0.89472824180426791
$ws.Range("K10").Value = 0.87042507387938794
$ws.Range("I11").Value = 0.83263079616127578
$ws.Range("J12").Value = 0.80697187391600766
$ws.Range("K12").Value = 0.86609611653811935
$ws.Range("K13").Value = 0.79200246376187411
$ws.Range("L13").Value = 0.9526924231815026
$ws.Range("M15").Value = 0.78144317153125742
$ws.Range("N16").Value = 0.85747353687954586
$ws.Range("O16").Value = 0.59555747477194299
$ws.Range("Q16").Value = 0.89574972586429347
$ws.Range("R16").Value = 0.97208133778272088
$ws.Range("S17").Value = 0.97526907281368191
$ws.Range("AR17").Value = 0.64632374459462416
$ws.Range("S18").Value = 0.85919255395153571
$ws.Range("AY19").Value = 0.85187171914823323
$ws.Range("BL19").Value = 0.68084193241752611
$ws.Range("BN19").Value = 0.9225524971424911
$ws.Range("T21").Value = 0.62025631776093926
$ws.Range("Y21").Value = 0.74359093767804074
$ws.Range("AV21").Value = 0.81137796133485074
$ws.Range("T22").Value = 0.90418432103836888
$ws.Range("A23").Value = 0.94522086499468427
$ws.Range("U23").Value = 0.93852669571943237
$ws.Range("V23").Value = 0.96528205170496784
$ws.Range("AA23").Value = 0.75824957172521268
$ws.Range("Y24").Value = 0.7653936091973369
$ws.Range("Z24").Value = 0.75014417337233763
$ws.Range("E25").Value = 0.95442261785471605
$ws.Range("W25").Value = 0.83404651215735526
$ws.Range("Z25").Value = 0.95738115582583172
$ws.Range("AL25").Value = 0.68711067589291819
$ws.Range("AA26").Value = 0.83661981092452198
$ws.Range("AB26").Value = 0.94090101461972753
$ws.Range("M27").Value = 0.94548715537342942
$ws.Range("AC27").Value = 0.95357887240642336
$ws.Range("AP27").Value = 0.87985265382232436
$ws.Range("AA28").Value = 0.89931910658634073
$ws.Range("AD28").Value = 0.78485260926738398
$ws.Range("AC30").Value = 0.95126623253291986
$ws.Range("AF30").Value = 0.9884422934116287
$ws.Range("AC31").Value = 0.81944591724818561
$ws.Range("AD31").Value = 0.92269406699968826
$ws.Range("AF31").Value = 0.93500124820146424
$ws.Range("AH32").Value = 0.99465682354594442
$ws.Range("BJ32").Value = 0.79709448379669945
$ws.Range("AE33").Value = 0.61987601709704421
$ws.Range("AF33").Value = 0.97167398918086356
$ws.Range("AI33").Value = 0.85530095778624116
$ws.Range("S34").Value = 0.88727440162987425
$ws.Range("AG34").Value = 0.81286474323849967
$ws.Range("AK35").Value = 0.89889086986481392
$ws.Range("N36").Value = 0.99320492841145569
$ws.Range("R36").Value = 0.82963477623975868
$ws.Range("AH36").Value = 0.74271336637853036
$ws.Range("AL36").Value = 0.8612791905957311
$ws.Range("AJ37").Value = 0.99399317122636166
$ws.Range("AL37").Value = 0.92480742591320064
$ws.Range("AM37").Value = 0.68687540061622465
$ws.Range("V39").Value = 0.8589793491286063
$ws.Range("BO39").Value = 0.6566445256938569
$ws.Range("AI40").Value = 0.90035743162332826
$ws.Range("AL40").Value = 0.94708048811733803
$ws.Range("AO40").Value = 0.86572118188119429
$ws.Range("H42").Value = 0.97190879989570234
$ws.Range("AN42").Value = 0.75417822237520138
$ws.Range("AO42").Value = 0.95721351066404203
$ws.Range("AQ42").Value = 0.99276506799122255
$ws.Range("E43").Value = 0.9768276958304839
$ws.Range("AO43").Value = 0.82368533078484263
$ws.Range("AR43").Value = 0.81271665876381283
$ws.Range("AT44").Value = 0.78297793727079024
$ws.Range("AQ45").Value = 0.95928661984226493
$ws.Range("AR45").Value = 0.94401931789353555
$ws.Range("AU45").Value = 0.51568317073173431
$ws.Range("AS46").Value = 0.68683341163600031
$ws.Range("AT47").Value = 0.78115341841879626
$ws.Range("AT48").Value = 0.89866559820890957
$ws.Range("AX48").Value = 0.76809041684249779
$ws.Range("AO49").Value = 0.96965534577586587
$ws.Range("AU49").Value = 0.91304738363832594
$ws.Range("AX49").Value = 0.99418314058246615
$ws.Range("AT50").Value = 0.93018989336171498
$ws.Range("W51").Value = 0.90780414382042518
$ws.Range("T52").Value = 0.81756491493354089
$ws.Range("AX52").Value = 0.93431181744159209
$ws.Range("BA52").Value = 0.88159781063717912
$ws.Range("BB53").Value = 0.82972068901905227
$ws.Range("BD54").Value = 0.89592698510498625
$ws.Range("BH54").Value = 0.89513294374160113
$ws.Range("X55").Value = 0.98405898444714512
$ws.Range("BA56").Value = 0.66288082015716721
$ws.Range("BF56").Value = 0.79954483735175663
$ws.Range("W57").Value = 0.52668978794064447
$ws.Range("AK57").Value = 0.79930414502999136
$ws.Range("BC57").Value = 0.91134483413096423
$ws.Range("BD57").Value = 0.86880039755450045
$ws.Range("BH57").Value = 0.89596694726789661
$ws.Range("BH58").Value = 0.93994268546399051
$ws.Range("BF59").Value = 0.74288961481423244
$ws.Range("L60").Value = 0.963239003841451
$ws.Range("BG60").Value = 0.90405466975615745
$ws.Range("BG61").Value = 0.73437299140670798
$ws.Range("BJ61").Value = 0.96130764084984166
$ws.Range("BK61").Value = 0.98076874045232953
$ws.Range("O62").Value = 0.97916980307927104
$ws.Range("AI63").Value = 0.75660942492567385
$ws.Range("BM63").Value = 0.65237221457811589
$ws.Range("BM64").Value = 0.80792664835156192
$ws.Range("BA65").Value = 0.67019341415239975
$ws.Range("BN65").Value = 0.6652554128113336
$ws.Range("BL66").Value = 0.82733943233518925
$ws.Range("BP66").Value = 0.78698640388705887
$ws.Range("A67").Value = 0.95558041402866434
$ws.Range("A68").Value = 0.70650295266369656
$ws.Range("B68").Value = 0.95586246769296102
$ws.Range("BO68").Value = 0.91006058438711579
